# Daily attendance processing - 2025-12-04 10:30:23
# For each row in column G ("Recorded By"), if the value is a comma-separated
# list of names ending in "System", rotate the list left by one position
# (move the first entry to the end) so "System" (or whichever entry was
# second-to-last) moves to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$col = 7  # column G

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "

    if ($parts.Count -ge 2 -and $parts[$parts.Count - 1] -eq "System") {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $newText = [string]::Join(", ", $rotated)
        $cell.Value2 = $newText
    }
}
